# Deltakerliste - KUNSTLØP: Oppvisning Bergen
# Remove the two cancelled ("Avmeldt") participants from the list:
#   - row 18: Mina Lam
#   - row 21: Oda Eilin Halkjelsvik-Sæbø
# Deleting the entire row shifts every row below it up by one, which is
# exactly what the target workbook shows (old row 19 -> new row 18, etc.,
# and the last two rows of the original sheet disappear).
#
# Delete the higher-numbered row first so the row number of the other
# target row does not shift before it is removed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(21).Delete()
$ws.Rows.Item(18).Delete()
